$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.794.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3807"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.260"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08232"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.550"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.474"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.657.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06979"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.808"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.790.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.084"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.231"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.842.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.940"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.084"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02837"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2532"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.155"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07140"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7094"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.348"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6579"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.340"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "

Write-Output "Applied all changes"